# Automatische test-sync: 2025-08-05 18:39:50
# Appends a new test-mail log row to the "Logs" sheet and bumps the
# matching rollup count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 37

$logs.Cells.Item($newRow, 1).Value = "Wil je dit even doorsturen?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #16: Wil je dit even doorsturen?"
$logs.Cells.Item($newRow, 4).Value = "Planning / Afspraak"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-05 18:39:28"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# "Planning / Afspraak" rollup count on the Dashboard goes from 19 to 20.
$dashboard.Range("B2").Value = 20

# Extend the conditional-formatting ranges to cover the newly added row.
$cols = "D", "G", "H", "I", "J"
foreach ($col in $cols) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "36")
    $newRange = $logs.Range("$col" + "2:" + "$col" + $newRow)
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
